$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '26.826.09'
Set-TextValue $ws.Range("E2") '  -1.66%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.857.23'
Set-TextValue $ws.Range("E3") '  -1.06%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.9998'
Set-TextValue $ws.Range("E4") '  -0.40%  '

# Row 5
Set-TextValue $ws.Range("D5") '304.93'
Set-TextValue $ws.Range("E5") '  -0.88%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.000'
Set-TextValue $ws.Range("E6") '  -0.25%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5067'
Set-TextValue $ws.Range("E7") '  -2.75%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3650'
Set-TextValue $ws.Range("E8") '  -2.88%  '

# Row 9
Set-TextValue $ws.Range("B9") 'Dogecoin'
Set-TextValue $ws.Range("C9") 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D9") '0.07150'
Set-TextValue $ws.Range("E9") '  -0.24%  '

# Row 10
Set-TextValue $ws.Range("B10") 'Polygon'
Set-TextValue $ws.Range("C10") 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D10") '0.8883'
Set-TextValue $ws.Range("E10") '  -0.01%  '

# Row 11
Set-TextValue $ws.Range("B11") 'Solana'
Set-TextValue $ws.Range("C11") 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range("D11") '20.66'
Set-TextValue $ws.Range("E11") '  -0.75%  '

# Row 12
Set-TextValue $ws.Range("B12") 'TRON'
Set-TextValue $ws.Range("C12") 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D12") '0.07527'
Set-TextValue $ws.Range("E12") '  -0.91%  '

# Row 13
Set-TextValue $ws.Range("B13") 'WrappedEther'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D13") '1.863.56'
Set-TextValue $ws.Range("E13") '  -1.72%  '

# Row 14
Set-TextValue $ws.Range("B14") 'Litecoin'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D14") '91.40'
Set-TextValue $ws.Range("E14") '  +2.16%  '

# Row 15
Set-TextValue $ws.Range("B15") 'Polkadot'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D15") '5.236'
Set-TextValue $ws.Range("E15") '  -2.00%  '

# Row 16
Set-TextValue $ws.Range("B16") 'BinanceUSD'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D16") '0.9998'
Set-TextValue $ws.Range("E16") '  -0.39%  '

# Row 17
Set-TextValue $ws.Range("B17") 'ShibaInu'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D17") '0.000008530'
Set-TextValue $ws.Range("E17") '  -0.54%  '

# Row 18
Set-TextValue $ws.Range("B18") 'Avalanche'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D18") '14.04'
Set-TextValue $ws.Range("E18") '  -1.19%  '

# Row 19
Set-TextValue $ws.Range("B19") 'Dai'
Set-TextValue $ws.Range("C19") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D19") '1.000'
Set-TextValue $ws.Range("E19") '  -0.26%  '

# Row 20
Set-TextValue $ws.Range("B20") 'WrappedBTC'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D20") '26.868.28'
Set-TextValue $ws.Range("E20") '  -1.67%  '

# Row 21
Set-TextValue $ws.Range("B21") 'Uniswap'
Set-TextValue $ws.Range("C21") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '5.010'
Set-TextValue $ws.Range("E21") '  -0.81%  '

# Row 22
Set-TextValue $ws.Range("B22") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C22") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D22") '2.094.63'
Set-TextValue $ws.Range("E22") '  -1.89%  '

# Row 23
Set-TextValue $ws.Range("B23") 'Cosmos'
Set-TextValue $ws.Range("C23") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D23") '10.27'
Set-TextValue $ws.Range("E23") '  -3.59%  '

# Row 24
Set-TextValue $ws.Range("B24") 'Chainlink'
Set-TextValue $ws.Range("C24") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D24") '6.436'
Set-TextValue $ws.Range("E24") '  -0.75%  '

# Row 25
Set-TextValue $ws.Range("B25") 'Toncoin'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D25") '1.822'
Set-TextValue $ws.Range("E25") '  -1.80%  '

# Row 26
Set-TextValue $ws.Range("B26") 'Monero'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D26") '146.55'
Set-TextValue $ws.Range("E26") '  -3.62%  '

# Row 27
Set-TextValue $ws.Range("B27") 'EthereumClassic'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D27") '17.85'
Set-TextValue $ws.Range("E27") '  -1.25%  '

# Row 28
Set-TextValue $ws.Range("B28") 'LidoDAOToken'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D28") '2.049'
Set-TextValue $ws.Range("E28") '  -5.84%  '

# Row 29
Set-TextValue $ws.Range("B29") 'BitcoinCash'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D29") '113.09'
Set-TextValue $ws.Range("E29") '  +0.13%  '

# Row 30
Set-TextValue $ws.Range("B30") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D30") '4.638'
Set-TextValue $ws.Range("E30") '  -2.51%  '

# Row 31
Set-TextValue $ws.Range("B31") 'Filecoin'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D31") '4.665'
Set-TextValue $ws.Range("E31") '  -1.03%  '

# Row 32
Set-TextValue $ws.Range("B32") 'Stellar'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D32") '0.09217'
Set-TextValue $ws.Range("E32") '  +1.85%  '

# Row 33
Set-TextValue $ws.Range("B33") 'Hedera'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D33") '0.05097'
Set-TextValue $ws.Range("E33") '  -1.71%  '

# Row 34
Set-TextValue $ws.Range("B34") 'HuobiToken'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D34") '3.069'
Set-TextValue $ws.Range("E34") '  -1.32%  '

# Row 35
Set-TextValue $ws.Range("B35") 'ARBITRUM'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D35") '1.147'
Set-TextValue $ws.Range("E35") '  -2.60%  '

# Row 36
Set-TextValue $ws.Range("B36") 'ImmutableX'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D36") '0.7311'
Set-TextValue $ws.Range("E36") '  -3.54%  '

# Row 37
Set-TextValue $ws.Range("B37") 'MXToken'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D37") '3.186'
Set-TextValue $ws.Range("E37") '  +4.40%  '

# Row 38
Set-TextValue $ws.Range("B38") 'VeChain'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D38") '0.02003'
Set-TextValue $ws.Range("E38") '  -2.11%  '

# Row 39
Set-TextValue $ws.Range("B39") 'RenderToken'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D39") '2.455'
Set-TextValue $ws.Range("E39") '  -3.16%  '

# Row 40
Set-TextValue $ws.Range("D40") '1.075'
Set-TextValue $ws.Range("E40") '  -0.79%  '

# Row 41
Set-TextValue $ws.Range("B41") 'TheSandbox'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D41") '0.5294'
Set-TextValue $ws.Range("E41") '  -3.08%  '

# Row 42
Set-TextValue $ws.Range("B42") 'Quant'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D42") '118.62'
Set-TextValue $ws.Range("E42") '  +2.36%  '

# Row 43
Set-TextValue $ws.Range("B43") 'FraxShare'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D43") '6.502'
Set-TextValue $ws.Range("E43") '  -2.70%  '

# Row 44
Set-TextValue $ws.Range("B44") 'Aptos'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D44") '8.413'
Set-TextValue $ws.Range("E44") '  -1.10%  '

# Row 45
Set-TextValue $ws.Range("B45") 'Algorand'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D45") '0.1471'
Set-TextValue $ws.Range("E45") '  -1.25%  '

# Row 46
Set-TextValue $ws.Range("B46") 'PaxDollar'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D46") '1.000'
Set-TextValue $ws.Range("E46") '  -0.21%  '

# Row 47
Set-TextValue $ws.Range("B47") 'Decentraland'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range("D47") '0.4626'
Set-TextValue $ws.Range("E47") '  -1.60%  '

# Row 48
Set-TextValue $ws.Range("B48") 'EnergySwap'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D48") '9.915'
Set-TextValue $ws.Range("E48") '  -3.12%  '

# Row 49
Set-TextValue $ws.Range("B49") 'Elrond'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range("D49") '37.27'
Set-TextValue $ws.Range("E49") '  +1.99%  '

# Row 50
Set-TextValue $ws.Range("B50") 'NEARProtocol'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D50") '1.558'
Set-TextValue $ws.Range("E50") '  -1.32%  '

# Row 51
Set-TextValue $ws.Range("B51") 'Aave'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D51") '62.88'
Set-TextValue $ws.Range("E51") '  -4.05%  '
